$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Student test 1"
$ws.Range("B3").Value = "Student test 1"

$ws.Range("C4").Select()
